# Update the "note" cell text (A5) from the old blue/green legend text to
# the new instructional text, and move the active-cell selection from D1
# down to A6 (below the table), matching the authored edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Green cells are entry locations. Any text outside these locations will be ignored"

$ws.Range("A6").Select() | Out-Null
